$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected OCR values for the frequency / dB table
$ws.Range("A1").Value = "50O OO0 MHz"
$ws.Range("B1").Value = "0 018 dB"
$ws.Range("C1").Value = "Not Valid **"
$ws.Range("D1").Value = "0 138 d8"

$ws.Range("A2").Value = "1.000 GHz"
$ws.Range("B2").Value = "-0.360 dB"
$ws.Range("C2").Value = "** Not Valid **"
$ws.Range("D2").Value = "-0.116 dB"

$ws.Range("A3").Value = "1.500 GHz"
$ws.Range("B3").Value = "-0.426 dB"
$ws.Range("C3").Value = "** Not Valid **"
$ws.Range("D3").Value = "-0.074 dB"

$ws.Range("A4").Value = "2.000 GHz"
$ws.Range("B4").Value = "-0.823 dB"
$ws.Range("C4").Value = "**Not Vaiid **"
$ws.Range("D4").Value = "-0.109 dB"

$ws.Range("A5").Value = "2.500 GHz"
$ws.Range("B5").Value = "-0.913 dB"
$ws.Range("C5").Value = "** Not Valid **"
$ws.Range("D5").Value = "-0.125 dB"

$ws.Range("A6").Value = "3.000 GHz"
$ws.Range("B6").Value = "-0.728 dB"
$ws.Range("C6").Value = "** Not Valid **"
$ws.Range("D6").Value = "-0.090 dB"

$ws.Range("A7").Value = "3.500 GHz"
$ws.Range("B7").Value = "-0.325 dB"
$ws.Range("C7").Value = "** Not Valid **"
$ws.Range("D7").Value = "-0.155 dB"

$ws.Range("A8").Value = "4.000 GHz"
$ws.Range("B8").Value = "-0.713 dB"
$ws.Range("C8").Value = "** Not Valid **"
$ws.Range("D8").Value = "-0.174 dB"

# Re-fit columns C/D now that the (longer) "** Not Valid **" text lives there
$ws.Columns.Item(3).ColumnWidth = 13.333333333333332
$ws.Columns.Item(4).ColumnWidth = 8.166666666666666
